# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that the file has been hung off for localization again:
#   - Status goes from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The "Latest Handoff File" / hyperlink display text is bumped to the
#     new handoff package name (b.63290e5768f688058c7b37413b0a5c26c308f864.*)
#   - The "Latest Handoff Datetime" / "Latest Handoff Date" timestamps are updated

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $readyForHandoff
$overview.Range("C3").Value = $readyForHandoff
$overview.Range("D3").Value = "2016-03-21 16:35:28"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $readyForHandoff
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-21 16:35:24"

$zhcnLinks = @()
foreach ($hl in $zhcn.Hyperlinks) {
    $zhcnLinks += $hl
}
foreach ($hl in $zhcnLinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyForHandoff
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-21 16:35:28"

$dedeLinks = @()
foreach ($hl in $dede.Hyperlinks) {
    $dedeLinks += $hl
}
foreach ($hl in $dedeLinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
